# Itération 2 - Finale
# Add the newly-collected participant rows (Nom / Prenom / grp) below the
# existing data table on "Feuil1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Fortier",   "Octave",   "A"),
  @("Brodeur",   "Alphonse", "A"),
  @("Bussiere",  "Matthieu", "A"),
  @("Mouet",     "Amaury",   "A"),
  @("Douffet",   "Estelle",  "A"),
  @("Franchet",  "Mayhew",   "A"),
  @("Compagnon", "Gregoire", "A"),
  @("Boulé",     "Magnolia", "A"),
  @("Corbin",    "Fiacre",   "A"),
  @("Dufresne",  "Lirienne", "A"),
  @("Françoise", "Boivin",   "B")
)

$startRow = 131
$row = $startRow
foreach ($entry in $data) {
  $ws.Cells.Item($row, 1).Value = $entry[0]
  $ws.Cells.Item($row, 2).Value = $entry[1]
  $ws.Cells.Item($row, 3).Value = $entry[2]
  $row++
}

# Leave the view roughly where the author left it: scrolled down near the
# bottom of the new data, with the last couple of touched cells selected.
$ws.Range("D136:D137").Select()
